$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$casesQuery = @"
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
 MATCH (f:file)-[*]->(c)
   WHERE c.race = "BLACK_OR_AFRICAN_AMERICAN"
RETURN DISTINCT
    c.case_id AS ``Case ID``,
     ct.clinical_trial_designation AS ``Trial Code``,
     a.arm_id AS Arm,
      a.arm_drug AS ``Arm Treatment``,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
"@.TrimEnd("`r","`n")

$filesQuery = @"
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
        WHERE c.race = "BLACK_OR_AFRICAN_AMERICAN"
WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS ``File Name``,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS ``File Format``,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS ``Trial Code``,
    a.arm_id AS Arm,
    c.case_id AS ``Case ID``
"@.TrimEnd("`r","`n")

$filesStatQuery = @"
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
      WHERE c.race = "BLACK_OR_AFRICAN_AMERICAN"
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
"@.TrimEnd("`r","`n")

$casesStatQuery = @"
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
         WHERE c.race = "BLACK_OR_AFRICAN_AMERICAN"
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
"@.TrimEnd("`r","`n")

# Write the new shared-string values in the same order the original authors
# introduced them (CasesTab already existed; FilesTab, the Cases query, the
# Files query, the Files stat-query and finally the Cases stat-query follow).
$ws.Range("A3").Value = "FilesTab"
$ws.Range("B2").Value = $casesQuery
$ws.Range("B3").Value = $filesQuery
$ws.Range("C3").Value = $filesStatQuery
$ws.Range("C2").Value = $casesStatQuery

$ws.Range("D3").Value = "TC03_Trials_Filter_Race-BlkAfrican_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC03_Trials_Filter_Race-BlkAfrican_WebData.xlsx"

$ws.Rows.Item(2).RowHeight = 195
$ws.Rows.Item(3).RowHeight = 409.5

# B3/C3 should wrap like B2/C2 already do
$ws.Range("B3:C3").WrapText = $true

# --- Selection moves to B2 ---
$ws.Range("B2").Select() | Out-Null
